# Auto-generated edit script: update Leviathan Profits workbook values
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 238398.72
$ws.Range("J17").Value = 238398.72
$ws.Range("L17").Value = 715196.16
$ws.Range("N17").Value = -715532.16
$ws.Range("H40").Value = 4834
$ws.Range("J40").Value = 4834
$ws.Range("L40").Value = 4834
$ws.Range("N40").Value = -5184
$ws.Range("H46").Value = 393867.78
$ws.Range("I46").Value = 973.26666
$ws.Range("J46").Value = 847207.6
$ws.Range("K46").Value = 2919.79998
$ws.Range("L46").Value = 2541622.8
$ws.Range("M46").Value = -2800.79998
$ws.Range("N46").Value = -2541860.8
$ws.Range("H60").Value = 393867.78
$ws.Range("I60").Value = 973.26666
$ws.Range("J60").Value = 847207.6
$ws.Range("K60").Value = 2919.79998
$ws.Range("L60").Value = 2541622.8
$ws.Range("M60").Value = -2435.79998
$ws.Range("N60").Value = -2542590.8
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H112").Value = 1430.6111
$ws.Range("I112").Value = 966.5
$ws.Range("J112").Value = 1563.2142
$ws.Range("K112").Value = 2899.5
$ws.Range("L112").Value = 4689.642599999999
$ws.Range("M112").Value = -1791.5
$ws.Range("N112").Value = -6905.642599999999
$ws.Range("H116").Value = 4999.8
$ws.Range("I116").Value = 4999
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 4999
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = -1557
$ws.Range("N116").Value = -11884
$ws.Range("H132").Value = 1250.8292
$ws.Range("I132").Value = 1207.125
$ws.Range("J132").Value = 2999
$ws.Range("K132").Value = 3621.375
$ws.Range("L132").Value = 8997
$ws.Range("M132").Value = -1091.375
$ws.Range("N132").Value = -14057
$ws.Range("H137").Value = 1714.2593
$ws.Range("I137").Value = 1664.619
$ws.Range("K137").Value = 4993.857
$ws.Range("M137").Value = -2443.857
$ws.Range("H140").Value = 188186.25
$ws.Range("J140").Value = 188186.25
$ws.Range("L140").Value = 188186.25
$ws.Range("N140").Value = -198546.25

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1936.9259
$ws.Range("I2").Value = 1932.8
$ws.Range("K2").Value = 1932.8
$ws.Range("M2").Value = -1819.8
$ws.Range("H32").Value = 24704.408
$ws.Range("I32").Value = 4319.5903
$ws.Range("K32").Value = 4319.5903
$ws.Range("M32").Value = -4032.5903
$ws.Range("H45").Value = 406247.62
$ws.Range("I45").Value = 595845.1
$ws.Range("K45").Value = 595845.1
$ws.Range("M45").Value = -595468.1
$ws.Range("H61").Value = 1631.5625
$ws.Range("I61").Value = 1631.5625
$ws.Range("K61").Value = 1631.5625
$ws.Range("M61").Value = -1419.5625
$ws.Range("H74").Value = 1553.52
$ws.Range("I74").Value = 1439.4736
$ws.Range("K74").Value = 1439.4736
$ws.Range("M74").Value = -565.4736
$ws.Range("H77").Value = 1553.52
$ws.Range("I77").Value = 1439.4736
$ws.Range("K77").Value = 7197.368
$ws.Range("M77").Value = -2829.368
$ws.Range("H80").Value = 20000
$ws.Range("J80").Value = 20000
$ws.Range("L80").Value = 20000
$ws.Range("N80").Value = -21996
$ws.Range("H83").Value = 20000
$ws.Range("J83").Value = 20000
$ws.Range("L83").Value = 60000
$ws.Range("N83").Value = -69984
$ws.Range("H116").Value = 1936.9259
$ws.Range("I116").Value = 1932.8
$ws.Range("K116").Value = 1932.8
$ws.Range("M116").Value = 361.2
$ws.Range("H122").Value = 1242.9166
$ws.Range("I122").Value = 1232.3265
$ws.Range("K122").Value = 3696.979499999999
$ws.Range("M122").Value = -1246.979499999999
$ws.Range("H132").Value = 1694
$ws.Range("I132").Value = 1587.85
$ws.Range("J132").Value = 2224.75
$ws.Range("K132").Value = 4763.549999999999
$ws.Range("L132").Value = 6674.25
$ws.Range("M132").Value = -2233.549999999999
$ws.Range("N132").Value = -11734.25
$ws.Range("H136").Value = 1631.5625
$ws.Range("I136").Value = 1631.5625
$ws.Range("K136").Value = 4894.6875
$ws.Range("M136").Value = -2344.6875

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1936.9259
$ws.Range("I3").Value = 1932.8
$ws.Range("K3").Value = 1932.8
$ws.Range("M3").Value = -1818.8
$ws.Range("H134").Value = 1641.3572
$ws.Range("I134").Value = 1431.5834
$ws.Range("K134").Value = 4294.7502
$ws.Range("M134").Value = -1759.7502
$ws.Range("H141").Value = 70340.375
$ws.Range("J141").Value = 71716.28999999999
$ws.Range("L141").Value = 71716.28999999999
$ws.Range("N141").Value = -82076.28999999999

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 5197.7407
$ws.Range("I132").Value = 4849.524
$ws.Range("J132").Value = 6416.5
$ws.Range("K132").Value = 14548.572
$ws.Range("L132").Value = 19249.5
$ws.Range("M132").Value = -12018.572
$ws.Range("N132").Value = -24309.5
$ws.Range("H134").Value = 2585.5557
$ws.Range("I134").Value = 2484.0466
$ws.Range("J134").Value = 2982.3635
$ws.Range("K134").Value = 7452.139800000001
$ws.Range("L134").Value = 8947.0905
$ws.Range("M134").Value = -4917.139800000001
$ws.Range("N134").Value = -14017.0905
$ws.Range("H137").Value = 80000
$ws.Range("J137").Value = 80000
$ws.Range("L137").Value = 80000
$ws.Range("N137").Value = -90200

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 347.72223
$ws.Range("I2").Value = 410
$ws.Range("J2").Value = 36.333332
$ws.Range("K2").Value = 2460
$ws.Range("L2").Value = 217.999992
$ws.Range("M2").Value = -2347
$ws.Range("N2").Value = -443.999992
$ws.Range("H38").Value = 53.57143
$ws.Range("I38").Value = 45.75
$ws.Range("K38").Value = 137.25
$ws.Range("M38").Value = 209.75
$ws.Range("H129").Value = 92923.73
$ws.Range("J129").Value = 4196
$ws.Range("L129").Value = 12588
$ws.Range("N129").Value = -22588
$ws.Range("H131").Value = 1541.4166
$ws.Range("J131").Value = 2032.5
$ws.Range("L131").Value = 6097.5
$ws.Range("N131").Value = -16177.5
$ws.Range("H137").Value = 9093005
$ws.Range("I137").Value = 10001307
$ws.Range("K137").Value = 30003921
$ws.Range("M137").Value = -29998821

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 22201
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H72").Value = 22201
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H126").Value = 3060.5
$ws.Range("I126").Value = 2642.4285
$ws.Range("K126").Value = 7927.2855
$ws.Range("M126").Value = -5457.2855
$ws.Range("H132").Value = 2439.077
$ws.Range("I132").Value = 2518.9092
$ws.Range("K132").Value = 7556.7276
$ws.Range("M132").Value = -5026.7276
$ws.Range("H136").Value = 48079.6
$ws.Range("J136").Value = 48079.6
$ws.Range("L136").Value = 144238.8
$ws.Range("N136").Value = -149338.8
$ws.Range("H138").Value = 81999.336
$ws.Range("J138").Value = 81999.336
$ws.Range("L138").Value = 81999.336
$ws.Range("N138").Value = -92279.336
$ws.Range("H139").Value = 99999
$ws.Range("J139").Value = 99999
$ws.Range("L139").Value = 99999
$ws.Range("N139").Value = -110279

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7004.923
$ws.Range("I7").Value = 7827.88
$ws.Range("J7").Value = 5535.357
$ws.Range("K7").Value = 7827.88
$ws.Range("L7").Value = 5535.357
$ws.Range("M7").Value = -7715.88
$ws.Range("N7").Value = -5759.357
$ws.Range("H46").Value = 60859.715
$ws.Range("I46").Value = 211761.5
$ws.Range("K46").Value = 211761.5
$ws.Range("M46").Value = -211573.5
$ws.Range("H93").Value = 18897.9
$ws.Range("I93").Value = 2435.0667
$ws.Range("K93").Value = 2435.0667
$ws.Range("M93").Value = -1187.0667
$ws.Range("H126").Value = 7004.923
$ws.Range("I126").Value = 7827.88
$ws.Range("J126").Value = 5535.357
$ws.Range("K126").Value = 23483.64
$ws.Range("L126").Value = 16606.071
$ws.Range("M126").Value = -21013.64
$ws.Range("N126").Value = -21546.071
$ws.Range("H132").Value = 3688.578
$ws.Range("I132").Value = 3237.0303
$ws.Range("J132").Value = 4930.3335
$ws.Range("K132").Value = 9711.090899999999
$ws.Range("L132").Value = 14791.0005
$ws.Range("M132").Value = -7181.090899999999
$ws.Range("N132").Value = -19851.0005

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2583
$ws.Range("I122").Value = 2199.6667
$ws.Range("K122").Value = 6599.000100000001
$ws.Range("M122").Value = -4149.000100000001
$ws.Range("H126").Value = 2971.2778
$ws.Range("I126").Value = 3127.4285
$ws.Range("K126").Value = 9382.2855
$ws.Range("M126").Value = -6912.2855
